$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in missing data for row 54 (四方坪站 on 45896)
$ws.Range("C54").Value = 11421.93
$ws.Range("D54").Value = 9865.16
$ws.Range("E54").Value = 3938.65
$ws.Range("F54").Value = 458

# Fill in missing data for row 55 (高岭站 on 45896)
$ws.Range("C55").Value = 5367.02
$ws.Range("D55").Value = 4589.21
$ws.Range("E55").Value = 1354.23
$ws.Range("F55").Value = 177

# Update the view: scroll position and selection
$ws.Range("I53").Select()
$excel.ActiveWindow.ScrollRow = 43
